$d = $word.ActiveDocument

function Set-RangeXml($rng, $bodyXml) {
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $bodyXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $rng.InsertXML($pkg)
}

# --- Step 1 (process from the end of the document backwards so earlier
#     paragraph indices/character offsets stay valid) ---

# 1a. Remove the trailing empty paragraph (last <w:p/>) by merging it away,
#     keeping paragraph 12 (with the bookmark) as the final paragraph.
$pCount = $d.Paragraphs.Count
$pLast = $d.Paragraphs($pCount)
$pPrev = $d.Paragraphs($pCount - 1)
$rngTrail = $d.Range($pPrev.Range.End - 1, $pLast.Range.End - 1)
$rngTrail.Delete()

# 1b. Remove the text run "Go ahead and get started." from what is now the
#     last paragraph, keeping the bookmarkStart/bookmarkEnd. (Use Text = ""
#     rather than Delete() -- Delete() on a range that empties a paragraph
#     can merge the paragraph away and drop its bookmarks.)
$pCount = $d.Paragraphs.Count
$pBookmark = $d.Paragraphs($pCount)
$rngGoAhead = $d.Range($pBookmark.Range.Start, $pBookmark.Range.Start + 26)
$rngGoAhead.Text = ""

# 1c. Delete the whole 3rd "section" entirely: the 2nd page-break paragraph,
#     the 3rd "Text Signatures" heading paragraph and the 3rd
#     "To take advantage..." paragraph (indices 9,10,11 of the original 13).
$pCount = $d.Paragraphs.Count
$pBreak2 = $d.Paragraphs($pCount - 3)
$pToTake3 = $d.Paragraphs($pCount - 1)
$rngSection3 = $d.Range($pBreak2.Range.Start, $pToTake3.Range.End)
$rngSection3.Delete()

# --- Step 2: split/merge runs & insert lastRenderedPageBreak ---

# 2a. Paragraph 7 (2nd "To take advantage...") -> split into two runs after
#     "format yo".
$marker = "format yo"
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute($marker) | Out-Null
$found = $d.Content.Find.Found

# locate paragraph 7 precisely
$p7 = $d.Paragraphs(7)
$full7 = $p7.Range.Text
$idx7 = $full7.IndexOf("format yo") + ("format yo").Length
$rng7 = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$xml7 = "<w:p><w:r><w:t xml:space='preserve'>" + $full7.Substring(0, $idx7) + "</w:t></w:r><w:r><w:t xml:space='preserve'>" + $full7.Substring($idx7) + "</w:t></w:r></w:p>"
Set-RangeXml $rng7 $xml7

# 2b. Paragraph 6 ("Text" + " Signatures") -> merge into a single run
#     "Text Signatures" and add <w:lastRenderedPageBreak/> before it.
$p6 = $d.Paragraphs(6)
$rng6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$xml6 = "<w:p><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:lastRenderedPageBreak/></w:r><w:r><w:t>Text Signatures</w:t></w:r></w:p>"
Set-RangeXml $rng6 $xml6

# 2c. Paragraph 3 (1st "To take advantage...") -> split into two runs after
#     "Emphasis and ".
$p3 = $d.Paragraphs(3)
$full3 = $p3.Range.Text
$idx3 = $full3.IndexOf("Emphasis and ") + ("Emphasis and ").Length
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$xml3 = "<w:p><w:r><w:t xml:space='preserve'>" + $full3.Substring(0, $idx3) + "</w:t></w:r><w:r><w:t xml:space='preserve'>" + $full3.Substring($idx3) + "</w:t></w:r></w:p>"
Set-RangeXml $rng3 $xml3

Write-Host "done"
